# Apply updated crypto price/volume data to Sheet1 (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(newPriceOrNull, newVolume)
$updates = @{
    2 = @("26.495.53", "  -3.27%  ")
    3 = @("1.805.99", "  -3.03%  ")
    4 = @("1.007", "  +0.33%  ")
    5 = @($null, "  +0.44%  ")
    6 = @("307.97", "  -2.30%  ")
    7 = @($null, "  -1.85%  ")
    8 = @("0.3657", "  -1.59%  ")
    9 = @("0.07104", "  -2.97%  ")
    10 = @("0.8731", "  -1.81%  ")
    11 = @("0.07787", "  -0.66%  ")
    12 = @("19.29", "  -3.83%  ")
    13 = @("1.833.27", "  -1.10%  ")
    14 = @("5.271", "  -2.37%  ")
    15 = @("6.336", "  -3.37%  ")
    16 = @("86.54", "  -5.77%  ")
    17 = @($null, "  +0.37%  ")
    18 = @("0.000008555", "  -4.79%  ")
    19 = @("1.007", "  +0.47%  ")
    20 = @("26.511.82", "  -3.25%  ")
    21 = @("14.20", "  -4.09%  ")
    22 = @("4.962", "  -3.41%  ")
    23 = @("2.062.78", "  -0.30%  ")
    24 = @("10.36", "  -1.97%  ")
    25 = @("1.977", "  +2.11%  ")
    26 = @("150.49", "  -0.93%  ")
    27 = @("17.81", "  -3.60%  ")
    28 = @("1.998", "  -2.94%  ")
    29 = @("113.22", "  -2.72%  ")
    30 = @("4.867", "  -4.69%  ")
    31 = @("0.08664", "  -2.10%  ")
    32 = @("3.117", "  -0.53%  ")
    33 = @("0.7283", "  -5.28%  ")
    34 = @("4.437", "  -1.80%  ")
    35 = @("1.114", "  -5.06%  ")
    36 = @("1.006", "  +0.65%  ")
    37 = @("2.543", "  -7.59%  ")
    38 = @("1.077", "  -0.40%  ")
    39 = @("0.01911", "  -2.57%  ")
    40 = @("0.05099", "  -2.81%  ")
    41 = @("2.865", "  -3.86%  ")
    42 = @("6.885", "  -2.76%  ")
    43 = @("0.4906", "  -4.64%  ")
    44 = @("0.1568", "  -4.66%  ")
    45 = @("8.139", "  -3.34%  ")
    46 = @($null, "  +0.55%  ")
    47 = @("0.4593", "  -4.58%  ")
    48 = @($null, "  -1.30%  ")
    49 = @("9.910", "  -4.86%  ")
    50 = @("1.580", "  -4.32%  ")
    51 = @("0.05995", "  -3.64%  ")
}

# Cells whose new price text would be auto-parsed as a number by Excel need to be
# forced to Text format first so the stored value matches the literal string.
$numericLooking = @(4, 6, 8, 9, 10, 11, 12, 14, 15, 16, 18, 19, 21, 22, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 47, 49, 50, 51)

foreach ($row in $numericLooking) {
    $ws.Range("D$row").NumberFormat = "@"
}

foreach ($row in $updates.Keys) {
    $priceText = $updates[$row][0]
    $volumeText = $updates[$row][1]
    if ($priceText -ne $null) {
        $ws.Range("D$row").Value = $priceText
    }
    $ws.Range("E$row").Value = $volumeText
}

# Restore default ("Normal") cell style so only the value/type changed, not formatting
foreach ($row in $numericLooking) {
    $ws.Range("D$row").Style = "Normal"
}
